$d = $word.ActiveDocument

# Move to the very end of the document content and insert a new paragraph
$end = $d.Content
$end.Collapse(0)            # wdCollapseEnd = 0
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Move(4, 1)              # wdCharacter = 1, move into the new paragraph

# Apply the same formatting used by the rest of the document
$end.Font.Name = "Times New Roman"
$end.Font.Bold = $true
$end.Font.Size = 12

$end.Text = "Im am studying EEE"
